$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "Price" cells are plain text in the source sheet (t="inlineStr").
# Excel COM auto-converts decimal-looking strings to numbers, so we force
# text interpretation with a leading single-quote (quote-prefix), matching
# how Excel itself preserves text such as "8.020" or "0.5351" verbatim.
$ws.Range("D2").Value = "'" + '26.349.05'
$ws.Range("E2").Value = '  +1.12%  '
$ws.Range("D3").Value = "'" + '1.667.09'
$ws.Range("E3").Value = '  +0.93%  '
$ws.Range("E4").Value = '  +0.93%  '
$ws.Range("E5").Value = '  +1.05%  '
$ws.Range("D6").Value = "'" + '0.5351'
$ws.Range("E6").Value = '  +1.70%  '
$ws.Range("E7").Value = '  +0.87%  '
$ws.Range("E8").Value = '  +2.49%  '
$ws.Range("E9").Value = '  +1.29%  '
$ws.Range("D10").Value = "'" + '20.92'
$ws.Range("E10").Value = '  +2.77%  '
$ws.Range("D11").Value = "'" + '0.07862'
$ws.Range("E11").Value = '  +0.85%  '
$ws.Range("D12").Value = "'" + '4.573'
$ws.Range("E12").Value = '  +1.47%  '
$ws.Range("D13").Value = "'" + '1.669.07'
$ws.Range("E13").Value = '  +1.27%  '
$ws.Range("D14").Value = "'" + '1.896.16'
$ws.Range("E14").Value = '  +0.95%  '
$ws.Range("D15").Value = "'" + '0.5540'
$ws.Range("E15").Value = '  +0.93%  '
$ws.Range("E16").Value = '  -0.13%  '
$ws.Range("D17").Value = "'" + '65.92'
$ws.Range("E17").Value = '  +0.64%  '
$ws.Range("D18").Value = "'" + '26.370.90'
$ws.Range("E19").Value = '  +0.92%  '
$ws.Range("D20").Value = "'" + '4.691'
$ws.Range("E20").Value = '  +2.72%  '
$ws.Range("D21").Value = "'" + '193.90'
$ws.Range("E21").Value = '  +1.64%  '
$ws.Range("E22").Value = '  +2.62%  '
$ws.Range("E23").Value = '  +0.31%  '
$ws.Range("E24").Value = '  +0.87%  '
$ws.Range("D25").Value = "'" + '146.48'
$ws.Range("E25").Value = '  +2.25%  '
$ws.Range("D26").Value = "'" + '0.1235'
$ws.Range("E26").Value = '  -0.14%  '
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("E28").Value = '  +0.44%  '
$ws.Range("D29").Value = "'" + '1.499'
$ws.Range("E29").Value = '  +4.69%  '
$ws.Range("D30").Value = "'" + '0.05878'
$ws.Range("E30").Value = '  +0.93%  '
$ws.Range("D31").Value = "'" + '1.284'
$ws.Range("E31").Value = '  +0.93%  '
$ws.Range("D32").Value = "'" + '3.623'
$ws.Range("E32").Value = '  +2.18%  '
$ws.Range("E33").Value = '  +0.78%  '
$ws.Range("E34").Value = '  +1.65%  '
$ws.Range("D35").Value = "'" + '0.9722'
$ws.Range("E35").Value = '  +2.67%  '
$ws.Range("D36").Value = "'" + '2.828'
$ws.Range("E36").Value = '  +1.76%  '
$ws.Range("D37").Value = "'" + '2.423'
$ws.Range("E37").Value = '  +0.56%  '
$ws.Range("D38").Value = "'" + '0.5845'
$ws.Range("E38").Value = '  +1.97%  '
$ws.Range("D39").Value = "'" + '0.01603'
$ws.Range("E39").Value = '  -0.47%  '
$ws.Range("D40").Value = "'" + '0.8637'
$ws.Range("E40").Value = '  +2.37%  '
$ws.Range("D41").Value = "'" + '1.066.41'
$ws.Range("E41").Value = '  +3.53%  '
$ws.Range("D42").Value = "'" + '5.842'
$ws.Range("E42").Value = '  +1.67%  '
$ws.Range("D43").Value = "'" + '1.011'
$ws.Range("E43").Value = '  +0.87%  '
$ws.Range("D44").Value = "'" + '104.79'
$ws.Range("E44").Value = '  +0.41%  '
$ws.Range("D45").Value = "'" + '1.807.56'
$ws.Range("E45").Value = '  +0.76%  '
$ws.Range("D46").Value = "'" + '57.90'
$ws.Range("E46").Value = '  +1.60%  '
$ws.Range("D47").Value = "'" + '0.0₈107'
$ws.Range("E47").Value = '  -4.97%  '
$ws.Range("E48").Value = '  +0.98%  '
$ws.Range("D49").Value = "'" + '0.4388'
$ws.Range("E49").Value = '  +1.52%  '
$ws.Range("D50").Value = "'" + '8.020'
$ws.Range("E50").Value = '  +2.48%  '
$ws.Range("E51").Value = '  +0.48%  '
